$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.202.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.001.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.80%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.996.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.67%  "
$ws.Range("E11").Value = "  +4.75%  "
$ws.Range("E12").Value = "  +5.37%  "
$ws.Range("E13").Value = "  +6.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.14%  "
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.180.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.500.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.003.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "454.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.86%  "
$ws.Range("E22").Value = "  +4.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.84%  "
$ws.Range("E25").Value = "  +15.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.15%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +18.26%  "
$ws.Range("E30").Value = "  +20.94%  "
$ws.Range("E31").Value = "  -4.91%  "
$ws.Range("E32").Value = "  +5.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.17%  "
$ws.Range("E34").Value = "  +4.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.994"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.22%  "
$ws.Range("E37").Value = "  +8.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +14.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.309"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.11%  "
$ws.Range("E42").Value = "  +7.56%  "
$ws.Range("E43").Value = "  +8.03%  "
$ws.Range("E44").Value = "  +4.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "395.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.66%  "
$ws.Range("E46").Value = "  +7.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.791.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.92%  "
$ws.Range("E51").Value = "  +4.86%  "
